# Fixed update to excel issue
# - Rename "Requested quantity" header to series-specific names on the
#   two existing sheets.
# - Add a new "PO Forecast" sheet (ds / PO_Forecast / yhat_lower / yhat_upper)
#   with the forecast series, formatted to match the existing sheets.

$wb = $excel.ActiveWorkbook

$weekly  = $wb.Worksheets.Item("Weekly Quantity")
$monthly = $wb.Worksheets.Item("Monthly Trend")

# Rename the "Requested quantity" column headers.
$weekly.Range("B1").Value  = "Weekly_PO_Qty"
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after "Monthly Trend".
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "PO Forecast"

# Header row.
$ws.Range("A1").Value = "ds"
$ws.Range("B1").Value = "PO_Forecast"
$ws.Range("C1").Value = "yhat_lower"
$ws.Range("D1").Value = "yhat_upper"

# Match the bold/bordered/centered header style used on the other sheets.
$weekly.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("D1").PasteSpecial(-4122)

# Forecast data rows: ds, PO_Forecast, yhat_lower, yhat_upper.
$data = @(
    @(45312.99999999999, 0, -344.8913957505893, 79.43642462206171),
    @(45361.99999999999, 0, -250.6053930148383, 179.8607496433505),
    @(45445.99999999999, 120, -91.85252092335469, 336.3122849661364),
    @(45452.99999999999, 133, -69.48116964857918, 337.7686593284902),
    @(45459.99999999999, 146, -69.94042437214122, 360.0482599983919),
    @(45466.99999999999, 159, -64.76651296190967, 363.7867663003507),
    @(45473.99999999999, 172, -27.03140211019504, 384.3669868376656),
    @(45480.99999999999, 185, -16.689439637796, 383.808309310126),
    @(45487.99999999999, 197, -9.083000635897639, 408.0734889637212),
    @(45494.99999999999, 210, 5.74375050023145, 424.0790683222779),
    @(45501.99999999999, 223, 14.24129708208288, 434.4657237931169),
    @(45508.99999999999, 236, 29.41068164052675, 441.326732151141),
    @(45515.99999999999, 249, 36.99919517461006, 444.6164490801217),
    @(45522.99999999999, 262, 38.81511723886322, 472.430285207048),
    @(45529.99999999999, 274, 65.09237655275901, 491.706650762917),
    @(45536.99999999999, 287, 77.36444204645929, 487.2825567815005),
    @(45543.99999999999, 300, 94.002214786078, 494.3916492086854),
    @(45550.99999999999, 313, 105.2928480819296, 531.8727944574343),
    @(45557.99999999999, 326, 123.8520269043571, 542.443192061745),
    @(45578.99999999999, 364, 145.6552915850334, 574.5929209652927),
    @(45585.99999999999, 377, 170.4020947717495, 584.0244957820146),
    @(45592.99999999999, 390, 172.1445564250488, 596.7196259413596),
    @(45599.99999999999, 403, 192.7230228741888, 612.5983421428283),
    @(45606.99999999999, 416, 210.4104467049697, 623.5305610183304),
    @(45613.99999999999, 428, 231.6992823292917, 644.4519300335496),
    @(45620.99999999999, 441, 234.4879009816692, 660.3934131120953),
    @(45627.99999999999, 454, 256.2012368352203, 672.3492434766398),
    @(45634.99999999999, 467, 278.5859057751361, 677.3647941033067),
    @(45641.99999999999, 480, 264.1697510143547, 678.809914980198),
    @(45648.99999999999, 492, 286.133296431012, 689.1069039644377),
    @(45655.99999999999, 505, 297.1381962080254, 735.0879147929671)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r += 1
}

# Match the date number format used for column A on the other sheets.
$weekly.Range("A2").Copy()
$ws.Range("A2:A32").PasteSpecial(-4122)

# Restore the original active sheet/selection (adding a sheet activates it).
$weekly.Activate() | Out-Null
$weekly.Range("A1").Select() | Out-Null
